# Weekly update: insert a new "Espinaca" price record as the new row 144
# (above the existing row that used to be row 144), pushing the existing
# rows 144:173 down to 145:174. All the constant/categorical columns for
# this sheet's data block (A, B, C, E, F, G, H, I, L, N, O, Q, R) keep the
# same values as the surrounding rows; only the weekly figures (D = Fecha,
# J = Volumen, K = Precio minimo, M = Precio promedio ponderado,
# P = Precio $/Kg) are new for this entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 144:173 down to 145:174, leaving a blank row 144.
$ws.Rows.Item(144).Insert()

# Fill in the new weekly record in row 144.
$ws.Cells.Item(144, 1).Value = 8
$ws.Cells.Item(144, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(144, 3).Value = "Coquimbo"
$ws.Cells.Item(144, 4).Value = 44504
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = 100112012
$ws.Cells.Item(144, 7).Value = "Espinaca"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 3000
$ws.Cells.Item(144, 11).Value = 400
$ws.Cells.Item(144, 12).Value = 500
$ws.Cells.Item(144, 13).Value = 450
$ws.Cells.Item(144, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(144, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(144, 16).Value = 900
$ws.Cells.Item(144, 17).Value = 0.5
$ws.Cells.Item(144, 18).Value = "Hortaliza"
